# The LOINC code "42348-3" row is being removed from the "Include #0"
# worksheet's Concept/Description table. That row currently lives at A12:B12.
# Deleting the entire row shifts every row below it up by one, which turns
# the old 14-row table (header + 13 concept rows incl. the trailing
# "System URI"/"http://loinc.org" row) into a 13-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Include #0")

$ws.Rows.Item(12).Delete()
